$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto values (price/volume) scraped on the new run date.
# Each cell is forced to Text format before assignment to avoid Excel
# auto-converting numeric-looking strings (e.g. '63.003.61') into numbers,
# then ClearFormats() strips the temporary style so no stray 's' attribute
# is left behind on the cell (matching the original inline-string layout).
$updates = @(
    @{ Cell = 'D2'; Value = '63.003.61' }
    @{ Cell = 'E2'; Value = '  -0.84%  ' }
    @{ Cell = 'D3'; Value = '3.049.88' }
    @{ Cell = 'E3'; Value = '  -1.24%  ' }
    @{ Cell = 'E4'; Value = '  -0.12%  ' }
    @{ Cell = 'D5'; Value = '582.55' }
    @{ Cell = 'E5'; Value = '  -1.52%  ' }
    @{ Cell = 'D6'; Value = '151.27' }
    @{ Cell = 'E6'; Value = '  -2.52%  ' }
    @{ Cell = 'E7'; Value = '  -0.07%  ' }
    @{ Cell = 'D8'; Value = '0.534' }
    @{ Cell = 'E8'; Value = '  -2.01%  ' }
    @{ Cell = 'D9'; Value = '3.050.42' }
    @{ Cell = 'E9'; Value = '  -0.93%  ' }
    @{ Cell = 'E10'; Value = '  -3.14%  ' }
    @{ Cell = 'D11'; Value = '5.82' }
    @{ Cell = 'E11'; Value = '  -0.48%  ' }
    @{ Cell = 'D12'; Value = '0.447' }
    @{ Cell = 'E12'; Value = '  -2.37%  ' }
    @{ Cell = 'E13'; Value = '  -3.13%  ' }
    @{ Cell = 'D14'; Value = '36.08' }
    @{ Cell = 'E14'; Value = '  -3.94%  ' }
    @{ Cell = 'E15'; Value = '  +1.84%  ' }
    @{ Cell = 'D16'; Value = '3.554.24' }
    @{ Cell = 'E16'; Value = '  -1.35%  ' }
    @{ Cell = 'D17'; Value = '7.13' }
    @{ Cell = 'E17'; Value = '  -0.95%  ' }
    @{ Cell = 'D18'; Value = '62.983.38' }
    @{ Cell = 'E18'; Value = '  -0.86%  ' }
    @{ Cell = 'D19'; Value = '3.054.68' }
    @{ Cell = 'E19'; Value = '  -1.09%  ' }
    @{ Cell = 'D20'; Value = '480.11' }
    @{ Cell = 'E20'; Value = '  +0.78%  ' }
    @{ Cell = 'D21'; Value = '14.30' }
    @{ Cell = 'E21'; Value = '  -2.60%  ' }
    @{ Cell = 'D22'; Value = '0.706' }
    @{ Cell = 'D23'; Value = '7.51' }
    @{ Cell = 'E23'; Value = '  -0.93%  ' }
    @{ Cell = 'D24'; Value = '2.39' }
    @{ Cell = 'E24'; Value = '  -0.86%  ' }
    @{ Cell = 'D25'; Value = '81.87' }
    @{ Cell = 'E25'; Value = '  +0.70%  ' }
    @{ Cell = 'D26'; Value = '12.65' }
    @{ Cell = 'E26'; Value = '  -2.17%  ' }
    @{ Cell = 'D27'; Value = '10.53' }
    @{ Cell = 'E27'; Value = '  +5.06%  ' }
    @{ Cell = 'D28'; Value = '0.999' }
    @{ Cell = 'E28'; Value = '  -0.02%  ' }
    @{ Cell = 'D29'; Value = '7.36' }
    @{ Cell = 'E29'; Value = '  -0.24%  ' }
    @{ Cell = 'E30'; Value = '  +0.03%  ' }
    @{ Cell = 'D31'; Value = '2.65' }
    @{ Cell = 'E31'; Value = '  -1.51%  ' }
    @{ Cell = 'D32'; Value = '2.20' }
    @{ Cell = 'E32'; Value = '  +0.39%  ' }
    @{ Cell = 'D33'; Value = '27.72' }
    @{ Cell = 'E33'; Value = '  +1.55%  ' }
    @{ Cell = 'D34'; Value = '0.110' }
    @{ Cell = 'E34'; Value = '  -2.88%  ' }
    @{ Cell = 'E35'; Value = '  +0.54%  ' }
    @{ Cell = 'D36'; Value = '0.0₃0810' }
    @{ Cell = 'E36'; Value = '  -4.51%  ' }
    @{ Cell = 'D37'; Value = '5.90' }
    @{ Cell = 'E37'; Value = '  -3.50%  ' }
    @{ Cell = 'D38'; Value = '2.19' }
    @{ Cell = 'E38'; Value = '  -1.60%  ' }
    @{ Cell = 'D39'; Value = '3.14' }
    @{ Cell = 'E39'; Value = '  -7.15%  ' }
    @{ Cell = 'D40'; Value = '9.18' }
    @{ Cell = 'E40'; Value = '  -1.82%  ' }
    @{ Cell = 'D41'; Value = '50.36' }
    @{ Cell = 'E41'; Value = '  -0.90%  ' }
    @{ Cell = 'D42'; Value = '426.59' }
    @{ Cell = 'E42'; Value = '  -3.92%  ' }
    @{ Cell = 'B43'; Value = 'TheGraph' }
    @{ Cell = 'C43'; Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt' }
    @{ Cell = 'D43'; Value = '0.286' }
    @{ Cell = 'E43'; Value = '  +0.41%  ' }
    @{ Cell = 'B44'; Value = 'Kaspa' }
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas' }
    @{ Cell = 'D44'; Value = '0.115' }
    @{ Cell = 'E44'; Value = '  +3.51%  ' }
    @{ Cell = 'B45'; Value = 'Maker' }
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr' }
    @{ Cell = 'D45'; Value = '2.843.54' }
    @{ Cell = 'E45'; Value = '  +1.48%  ' }
    @{ Cell = 'B46'; Value = 'VeChain' }
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet' }
    @{ Cell = 'D46'; Value = '0.0361' }
    @{ Cell = 'E46'; Value = '  -0.59%  ' }
    @{ Cell = 'D47'; Value = '37.86' }
    @{ Cell = 'E47'; Value = '  -5.63%  ' }
    @{ Cell = 'D48'; Value = '127.10' }
    @{ Cell = 'E48'; Value = '  -3.22%  ' }
    @{ Cell = 'D50'; Value = '25.13' }
    @{ Cell = 'E50'; Value = '  -1.56%  ' }
    @{ Cell = 'D51'; Value = '0.109' }
    @{ Cell = 'E51'; Value = '  -1.12%  ' }
)

foreach ($u in $updates) {
    $rng = $ws.Range($u.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $u.Value
    $rng.ClearFormats()
}
